$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new part row (row 45): Micro USB B Male to USB B Female ---
$ws.Range("A45").Value = "Micro USB B Male to USB B Female"
$ws.Range("A45").Style = "Hyperlink"

# Copy the numeric/text formatting from row 44 (C:G) down into row 45
# so the new row matches the look of the rest of the table.
$ws.Range("C44:G44").Copy()
$ws.Range("C45:G45").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 7.58
$ws.Range("G45").Value = "Amazon"

# Hyperlink for the new part name cell
$ws.Hyperlinks.Add($ws.Range("A45"), "https://www.amazon.com/gp/product/B01N9RG5L8/ref=ppx_yo_dt_b_asin_title_o00_s00?ie=UTF8&psc=1")

# --- Update the view: scroll position / selection moved while editing ---
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B42").Select()
